$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").Value = "kontrola a oponentura 3.iterace ŠIP"
$ws.Range("B45").Value = 1.5

$ws.Range("A46").Select()
